$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on all D (Price) and E (Volume) cells being updated so that
# numeric-looking strings (e.g. "1.00", "467.30") are preserved verbatim as text
# instead of being coerced into numbers by COM Value assignment.
$cellRefs = @(
    "D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "E7", "E8", "E9", "D10", "E10", "E12", "D13", "E13", "E14", "D15", "E15", "D16", "E16", "E17", "D18", "E18", "E19", "D20", "E20", "E21", "E22", "D24", "E24", "E26", "D27", "E27", "E28", "E29", "D31", "E31", "E32", "E33", "D34", "E34", "E35", "D36", "E36", "D37", "E37", "D38", "E38", "D40", "E40", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "D48", "E48", "E49", "D50", "E50", "D51", "E51"
)
foreach ($ref in $cellRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = "69.252.18"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "3.673.51"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "675.25"
$ws.Range("E5").Value = "  -0.95%  "
$ws.Range("D6").Value = "157.67"
$ws.Range("E6").Value = "  -3.14%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("E8").Value = "  -1.42%  "
$ws.Range("E9").Value = "  -1.78%  "
$ws.Range("D10").Value = "6.96"
$ws.Range("E10").Value = "  -5.87%  "
$ws.Range("E12").Value = "  -3.41%  "
$ws.Range("D13").Value = "4.292.18"
$ws.Range("E13").Value = "  -0.49%  "
$ws.Range("E14").Value = "  -4.07%  "
$ws.Range("D15").Value = "3.674.35"
$ws.Range("E15").Value = "  -0.35%  "
$ws.Range("D16").Value = "69.192.10"
$ws.Range("E16").Value = "  -0.49%  "
$ws.Range("E17").Value = "  +1.51%  "
$ws.Range("D18").Value = "16.03"
$ws.Range("E18").Value = "  -1.60%  "
$ws.Range("E19").Value = "  -3.16%  "
$ws.Range("D20").Value = "467.30"
$ws.Range("E20").Value = "  -3.26%  "
$ws.Range("E21").Value = "  +0.53%  "
$ws.Range("E22").Value = "  -2.96%  "
$ws.Range("D24").Value = "3.816.81"
$ws.Range("E24").Value = "  -0.48%  "
$ws.Range("E26").Value = "  -7.21%  "
$ws.Range("D27").Value = "10.91"
$ws.Range("E27").Value = "  -5.11%  "
$ws.Range("E28").Value = "  -5.45%  "
$ws.Range("E29").Value = "  -1.71%  "
$ws.Range("D31").Value = "6.62"
$ws.Range("E31").Value = "  -4.09%  "
$ws.Range("E32").Value = "  -0.08%  "
$ws.Range("E33").Value = "  -5.17%  "
$ws.Range("D34").Value = "26.84"
$ws.Range("E34").Value = "  -1.19%  "
$ws.Range("E35").Value = "  +0.16%  "
$ws.Range("D36").Value = "0.161"
$ws.Range("E36").Value = "  -4.53%  "
$ws.Range("D37").Value = "8.12"
$ws.Range("E37").Value = "  -4.79%  "
$ws.Range("D38").Value = "6.21"
$ws.Range("E38").Value = "  -2.04%  "
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("E41").Value = "  -1.23%  "
$ws.Range("D42").Value = "175.05"
$ws.Range("E42").Value = "  +8.40%  "
$ws.Range("D43").Value = "0.0899"
$ws.Range("E43").Value = "  -4.48%  "
$ws.Range("D44").Value = "0.940"
$ws.Range("E44").Value = "  -1.68%  "
$ws.Range("D45").Value = "47.56"
$ws.Range("E45").Value = "  -1.66%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "27.99"
$ws.Range("E46").Value = "  -7.32%  "
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").Value = "2.69"
$ws.Range("E47").Value = "  -5.76%  "
$ws.Range("D48").Value = "0.000275"
$ws.Range("E48").Value = "  -4.81%  "
$ws.Range("E49").Value = "  -5.12%  "
$ws.Range("D50").Value = "1.07"
$ws.Range("E50").Value = "  -3.82%  "
$ws.Range("D51").Value = "7.77"
$ws.Range("E51").Value = "  -3.32%  "
